$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same summary data and need
# the same updates to column F ("想去人数" / want-to-go count).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2184
    $ws.Range("F4").Value = 1564
    $ws.Range("F5").Value = 7330
}
